# Updates cryptos list values (Price/Volume columns) and two name/link swaps
# per the "Updated cryptos list ... with GitHub Actions" commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.517.63'
$ws.Range("E2").Value = '  -1.25%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.572.85'
$ws.Range("E3").Value = '  -2.45%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '516.20'
$ws.Range("E5").Value = '  -2.69%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '138.44'
$ws.Range("E6").Value = '  -4.58%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.997'
$ws.Range("E7").Value = '  -0.10%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.563'
$ws.Range("E8").Value = '  -1.21%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.588.80'
$ws.Range("E9").Value = '  -2.26%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.48'
$ws.Range("E10").Value = '  -2.39%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.0999'
$ws.Range("E11").Value = '  -3.10%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.330'
$ws.Range("E12").Value = '  -1.17%  '
$ws.Range("E13").Value = '  +1.09%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.025.92'
$ws.Range("E14").Value = '  -2.42%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.473.38'
$ws.Range("E15").Value = '  -1.28%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.26'
$ws.Range("E16").Value = '  -1.74%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.565.44'
$ws.Range("E17").Value = '  -1.92%  '
$ws.Range("E18").Value = '  -3.07%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '336.60'
$ws.Range("E19").Value = '  -1.88%  '
$ws.Range("E20").Value = '  -2.83%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.09'
$ws.Range("E21").Value = '  -3.89%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.42'
$ws.Range("E22").Value = '  +1.49%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.11%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '65.84'
$ws.Range("E24").Value = '  -1.10%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.166'
$ws.Range("E25").Value = '  -0.73%  '
$ws.Range("E26").Value = '  -2.71%  '
$ws.Range("E27").Value = '  -0.24%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.99'
$ws.Range("E28").Value = '  -1.84%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0708'
$ws.Range("E30").Value = '  -10.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.92'
$ws.Range("E31").Value = '  -5.74%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '18.63'
$ws.Range("E32").Value = '  -1.57%  '
$ws.Range("E33").Value = '  -3.45%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.13'
$ws.Range("E34").Value = '  -1.20%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.92'
$ws.Range("E35").Value = '  -4.74%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.12'
$ws.Range("E36").Value = '  -3.48%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.25'
$ws.Range("E37").Value = '  -0.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.835'
$ws.Range("E38").Value = '  -0.23%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.43'
$ws.Range("E39").Value = '  -1.59%  '
$ws.Range("B40").Value = 'Fetch.AI'
$ws.Range("C40").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.820'
$ws.Range("E40").Value = '  -2.30%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.51'
$ws.Range("E41").Value = '  -2.42%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.997'
$ws.Range("E42").Value = '  +0.01%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '271.65'
$ws.Range("E43").Value = '  +1.01%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.72'
$ws.Range("E44").Value = '  +0.07%  '
$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0945'
$ws.Range("E45").Value = '  -2.87%  '
$ws.Range("B46").Value = 'Mantle'
$ws.Range("C46").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.584'
$ws.Range("E46").Value = '  -2.02%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0515'
$ws.Range("E47").Value = '  -3.01%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.968.01'
$ws.Range("E48").Value = '  -3.00%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '18.29'
$ws.Range("E49").Value = '  -3.96%  '
$ws.Range("E50").Value = '  -4.35%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.41'
$ws.Range("E51").Value = '  -5.90%  '
